$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D6").Value = "wrong format `$100"
$ws.Range("C6").Value = "wrong unit `$100"
$ws.Range("E6").Value = "different format `$100 but save as string"
$ws.Range("A4").Value = "blah blah"

$ws.Range("B17").Select()

$ws.Columns("E").ColumnWidth = 30
